# Generate Report for Handoff
#
# The "53659c0d-1d70-41d4-b83d-d99d379c43e0.md" file has finished
# translation for both locales and is now ready to be handed off.
# Update its Status + Latest Handoff Datetime on the Overview sheet and
# on each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the 53659c0d-...md file ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-21 10:17:16"

# --- zh-cn sheet: row 3 is the 53659c0d-...md file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-21 10:17:13"

# --- de-de sheet: row 3 is the 53659c0d-...md file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-21 10:17:16"
